# "Datos arreglados y añadidos subrutina form_kv"
# Fix the data in the "Elementos" sheet (element end-node numbering got
# shuffled), clean up the stray fill style applied to the "Nodos" FX column,
# and restore the UI selection state (active sheet back to "Elementos").

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Elementos")
$ws2 = $wb.Worksheets.Item("Nodos")
$ws3 = $wb.Worksheets.Item("Datos")

# --- Fix element connectivity values on "Elementos" ---------------------
$ws1.Range("C2").Value = 2
$ws1.Range("C3").Value = 1

$ws1.Range("B4").Value = 2
$ws1.Range("C4").Value = 1

$ws1.Range("B5").Value = 2
$ws1.Range("C5").Value = 4

$ws1.Range("B7").Value = 1
$ws1.Range("C7").Value = 3

$ws1.Range("B8").Value = 4
$ws1.Range("C8").Value = 3

$ws1.Range("A9").Value = 1
$ws1.Range("A10").Value = 1
$ws1.Range("A11").Value = 1

$ws1.Range("A12").Value = 1
$ws1.Range("B12").Value = 6
$ws1.Range("C12").Value = 5

$ws1.Range("A13").Value = 1
$ws1.Range("B13").Value = 6
$ws1.Range("C13").Value = 8

$ws1.Range("B15").Value = 5
$ws1.Range("C15").Value = 7

$ws1.Range("A17").Value = 1
$ws1.Range("A18").Value = 1
$ws1.Range("A19").Value = 1
$ws1.Range("A20").Value = 1

$ws1.Range("B21").Value = 10
$ws1.Range("B22").Value = 9

# --- Clean up the stray "applyFill" style on the FX column of "Nodos" ---
# (re-applying the existing thin border collapses the cell style back onto
# the plain bordered style, dropping the unused fill-flagged style)
$ws2.Range("F1:F13").Borders.LineStyle = 1

# --- Restore selections / active sheet -----------------------------------
[void]$ws3.Range("A3").Select()

[void]$ws2.Range("E3").Select()

[void]$ws1.Activate()
[void]$ws1.Range("L11").Select()
